$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row from hunk 0
$ws.Range("H55").Value = 86.666664
$ws.Range("I55").Value = 90.5
$ws.Range("J55").Value = 79
$ws.Range("K55").Value = 90.5
$ws.Range("L55").Value = 79
$ws.Range("M55").Value = 123.5
$ws.Range("N55").Value = -507
# row from hunk 1
$ws.Range("H93").Value = 88887
$ws.Range("J93").Value = 88887
$ws.Range("L93").Value = 88887
$ws.Range("N93").Value = -93879
# row from hunk 2
$ws.Range("H98").Value = 2334.111
$ws.Range("I98").Value = 2499.9092
$ws.Range("J98").Value = 2073.5715
$ws.Range("K98").Value = 2499.9092
$ws.Range("L98").Value = 2073.5715
$ws.Range("M98").Value = -1001.9092
$ws.Range("N98").Value = -5069.5715
# row from hunk 3
$ws.Range("H122").Value = 2334.111
$ws.Range("I122").Value = 2499.9092
$ws.Range("J122").Value = 2073.5715
$ws.Range("K122").Value = 7499.7276
$ws.Range("L122").Value = 6220.7145
$ws.Range("M122").Value = -5049.7276
$ws.Range("N122").Value = -11120.7145
# row from hunk 4
$ws.Range("H131").Value = 2574.6667
$ws.Range("J131").Value = 3222.3333
$ws.Range("L131").Value = 9666.999899999999
$ws.Range("N131").Value = -19746.9999
# row from hunk 5
$ws.Range("H132").Value = 1085.52
$ws.Range("I132").Value = 1066.9131
$ws.Range("K132").Value = 3200.7393
$ws.Range("M132").Value = -670.7393000000002
# row from hunk 6
$ws.Range("H135").Value = 620.4
$ws.Range("I135").Value = 451.5
$ws.Range("J135").Value = 873.75
$ws.Range("K135").Value = 4063.5
$ws.Range("L135").Value = 7863.75
$ws.Range("M135").Value = -1528.5
$ws.Range("N135").Value = -12933.75
# row from hunk 9
$ws.Range("H95").Value = 28000
$ws.Range("J95").Value = 28000
$ws.Range("L95").Value = 28000
$ws.Range("N95").Value = -33492

$ws = $wb.Worksheets.Item("ARM")
# row from hunk 7
$ws.Range("H32").Value = 2426.1875
$ws.Range("I32").Value = 1881.4117
$ws.Range("K32").Value = 1881.4117
$ws.Range("M32").Value = -1594.4117
# row from hunk 8
$ws.Range("H61").Value = 3163.8386
$ws.Range("I61").Value = 2462.8948
$ws.Range("K61").Value = 2462.8948
$ws.Range("M61").Value = -2250.8948
# row from hunk 10
$ws.Range("H136").Value = 3163.8386
$ws.Range("I136").Value = 2462.8948
$ws.Range("K136").Value = 7388.6844
$ws.Range("M136").Value = -4838.6844
# row from hunk 31
$ws.Range("H135").Value = 37000
$ws.Range("J135").Value = 37000
$ws.Range("L135").Value = 37000
$ws.Range("N135").Value = -47140

$ws = $wb.Worksheets.Item("CRP")
# row from hunk 11
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471
# row from hunk 12
$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181
# row from hunk 13
$ws.Range("H31").Value = 1323.4474
$ws.Range("I31").Value = 841.75
$ws.Range("J31").Value = 1858.6666
$ws.Range("K31").Value = 841.75
$ws.Range("L31").Value = 1858.6666
$ws.Range("M31").Value = -546.75
$ws.Range("N31").Value = -2448.6666
# row from hunk 14
$ws.Range("H34").Value = 1323.4474
$ws.Range("I34").Value = 841.75
$ws.Range("J34").Value = 1858.6666
$ws.Range("K34").Value = 841.75
$ws.Range("L34").Value = 1858.6666
$ws.Range("M34").Value = -639.75
$ws.Range("N34").Value = -2262.6666
# row from hunk 15
$ws.Range("H58").Value = 3624510.2
$ws.Range("I58").Value = 8696934
$ws.Range("K58").Value = 8696934
$ws.Range("M58").Value = -8696731
# row from hunk 16
$ws.Range("H92").Value = 39999
$ws.Range("J92").Value = 39999
$ws.Range("L92").Value = 39999
$ws.Range("N92").Value = -44991
# row from hunk 17
$ws.Range("H93").Value = 6666.3335
$ws.Range("I93").Value = 6666.3335
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 6666.3335
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -4794.3335
$ws.Range("N93").ClearContents()
# row from hunk 18
$ws.Range("H94").Value = 695.3158
$ws.Range("I94").Value = 571.8
$ws.Range("J94").Value = 832.55554
$ws.Range("K94").Value = 571.8
$ws.Range("L94").Value = 832.55554
$ws.Range("M94").Value = -120.8
$ws.Range("N94").Value = -1734.55554
# row from hunk 19
$ws.Range("H99").Value = 2599.8333
$ws.Range("I99").Value = 2419.8
$ws.Range("K99").Value = 2419.8
$ws.Range("M99").Value = -921.8000000000002
# row from hunk 20
$ws.Range("H126").Value = 2599.8333
$ws.Range("I126").Value = 2419.8
$ws.Range("K126").Value = 7259.400000000001
$ws.Range("M126").Value = -4789.400000000001
# row from hunk 21
$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959
# row from hunk 22
$ws.Range("H134").Value = 2345.5557
$ws.Range("I134").Value = 1077.75
$ws.Range("J134").Value = 3359.8
$ws.Range("K134").Value = 3233.25
$ws.Range("L134").Value = 10079.4
$ws.Range("M134").Value = -698.25
$ws.Range("N134").Value = -15149.4
# row from hunk 23
$ws.Range("H136").Value = 3624510.2
$ws.Range("I136").Value = 8696934
$ws.Range("K136").Value = 26090802
$ws.Range("M136").Value = -26088252

$ws = $wb.Worksheets.Item("CUL")
# row from hunk 24
$ws.Range("H131").Value = 7475771.5
$ws.Range("I131").Value = 166667140
$ws.Range("K131").Value = 500001420
$ws.Range("M131").Value = -499996380
# row from hunk 25
$ws.Range("H132").Value = 997.5217
$ws.Range("I132").Value = 862
$ws.Range("J132").Value = 1035.1666
$ws.Range("K132").Value = 7758
$ws.Range("L132").Value = 9316.499400000001
$ws.Range("M132").Value = -5228
$ws.Range("N132").Value = -14376.4994

$ws = $wb.Worksheets.Item("GSM")
# row from hunk 26
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# row from hunk 27
$ws.Range("H22").Value = 3897.7778
$ws.Range("J22").Value = 5956
$ws.Range("L22").Value = 5956
$ws.Range("N22").Value = -6546
# row from hunk 28
$ws.Range("H27").Value = 3897.7778
$ws.Range("J27").Value = 5956
$ws.Range("L27").Value = 5956
$ws.Range("N27").Value = -6170
# row from hunk 29
$ws.Range("H61").Value = 2787.65
$ws.Range("I61").Value = 2013
$ws.Range("J61").Value = 4226.2856
$ws.Range("K61").Value = 2013
$ws.Range("L61").Value = 4226.2856
$ws.Range("M61").Value = -1811
$ws.Range("N61").Value = -4630.2856
# row from hunk 30
$ws.Range("H113").Value = 2787.65
$ws.Range("I113").Value = 2013
$ws.Range("J113").Value = 4226.2856
$ws.Range("K113").Value = 2013
$ws.Range("L113").Value = 4226.2856
$ws.Range("M113").Value = 157
$ws.Range("N113").Value = -8566.285599999999

$ws = $wb.Worksheets.Item("WVR")
# row from hunk 32
$ws.Range("H122").Value = 50544.152
$ws.Range("I122").Value = 62230.906
$ws.Range("K122").Value = 186692.718
$ws.Range("M122").Value = -184242.718
# row from hunk 33
$ws.Range("H126").Value = 4534.3335
$ws.Range("J126").Value = 6995
$ws.Range("L126").Value = 20985
$ws.Range("N126").Value = -25925
